$d = $word.ActiveDocument

function Add-RawParagraphXml($innerBodyXml) {
    $endPos = $d.Content.End
    $r = $d.Range($endPos, $endPos)
    $pkg = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $innerBodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($pkg)
}

# --- Change 1: paragraph 3 (end of RF1 section) ---
$old1 = 'Gesiane digita sua senha e aperta em confirmar. Após o sistema verificar a senha, ele pede para que ela retire o seu cartão de débito e aguarde a impressão do QR Code e logo depois que o mesmo imprime, ela usa o bilhete para passar pela catraca da estação. '
$new1 = 'Gesiane digita sua senha e aperta em confirmar. Após o sistema verificar a senha, ele pede para que ela retire o seu cartão de débito e aguarde a impressão do QR Code. Depois disso, o sistema exibe uma mensagem de que a transação foi efetuada com sucesso e volta à sua tela inicial depois de alguns segundos.'
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# --- Change 2: paragraph 4 (end of RF1 section) ---
$old2 = 'E assim a mesma pode aguardar o trem chegar para que possa seguir para seu destino.'
$new2 = 'E logo depois que o mesmo imprime, ela usa o bilhete para passar pela catraca da estação. E assim a mesma pode aguardar o trem chegar para que possa seguir para seu destino.'
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# --- Change 3: paragraph 9 (end of RF2 section, accent fix "a sua" -> "a sua" with grave) ---
$old3 = 'Prosseguindo, Gesiane escolhe carregar seu cartão com R$50,00 e o sistema pede para que ela insira as notas até completar o valor escolhido. Após isso, o sistema exibe uma mensagem de que o valor inserido na máquina foi aceito caso esteja certo e deixa Gesiane escolher se quer o recibo impresso ou não. Ela seleciona que não e o sistema pede para que a mesma retire seu cartão top. Depois disso, o sistema exibe uma mensagem de que a transação foi efetuada com sucesso e volta a sua tela inicial depois de alguns segundos.'
$new3 = 'Prosseguindo, Gesiane escolhe carregar seu cartão com R$50,00 e o sistema pede para que ela insira as notas até completar o valor escolhido. Após isso, o sistema exibe uma mensagem de que o valor inserido na máquina foi aceito caso esteja certo e deixa Gesiane escolher se quer o recibo impresso ou não. Ela seleciona que não e o sistema pede para que a mesma retire seu cartão top. Depois disso, o sistema exibe uma mensagem de que a transação foi efetuada com sucesso e volta à sua tela inicial depois de alguns segundos.'
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2) | Out-Null

# --- Change 4: append new RF3 section (6 new paragraphs) at end of document ---
Add-RawParagraphXml('<w:p><w:pPr><w:ind w:left="0" w:firstLine="0"/><w:jc w:val="both"/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p><w:p><w:pPr><w:ind w:left="0" w:firstLine="0"/><w:jc w:val="both"/></w:pPr><w:r><w:rPr><w:b w:val="1"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">RF3 - </w:t></w:r><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Gesiane agora, tem seu cartão de débito funcionando novamente e como sempre fez, agora pretende carregar seu cartão top com a forma de pagamento de débito.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="0" w:firstLine="0"/><w:jc w:val="both"/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:tab/><w:t xml:space="preserve">Gesiane chega na máquina para recarregar seu cartão e a mesma seleciona a opção de recarregar seu cartão top. Após selecionar a opção, a máquina pede para que ela escolha a forma de pagamento. Ela seleciona a opção para pagar com o cartão de débito e logo em seguida o sistema pede à ela que a mesma insira o seu cartão top no local indicado. Depois de inserir o cartão, na tela é mostrado uma tela para a seleção do tipo de recarga do cartão top: comum ou escolar. Na tela é mostrado o valor limite que pode ser carregado no cartão. E esse valor é de R$3.000,00.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:tab/><w:t xml:space="preserve">Gesiane seleciona a opção “Comum” e prossegue com a transação. O sistema pede para ela inserir seu cartão de débito e depois dessa ação, solicita a senha do cartão para Gesiane depois de alguns segundos.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:tab/><w:t xml:space="preserve">Gesiane digita sua senha e aperta em confirmar. Após o sistema verificar a senha, ele pede para que ela retire o seu cartão de débito e aguarde a impressão do QR Code. Depois disso, o sistema exibe uma mensagem de que a transação foi efetuada com sucesso e volta à sua tela inicial depois de alguns segundos.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="720"/><w:jc w:val="both"/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">E logo depois que o mesmo imprime, ela usa o bilhete para passar pela catraca da estação. E assim a mesma pode aguardar o trem chegar para que possa seguir para seu destino. </w:t></w:r></w:p>')

